$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (Volume number + report week dates)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# ---------------------------------------------------------------------
# Helper: set a numeric cell making sure the style keeps the "numeric"
# look for its column (re-applying NumberFormat nudges the engine to
# reuse the existing numeric style instead of leaving a stale
# text-flavoured style behind after a Value assignment).
# ---------------------------------------------------------------------

# Row 15 (Rape): C/D/E switch from numbers to the "no data" text markers,
# mirroring the formatting already used on row 14 (style 13, shared
# strings "0" / "***.*"). Copying from row 14 preserves both the shared
# string reuse and the exact style index.
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("G15").Value = 2
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 16 (Robbery)
$ws.Range("D16").Value = 9
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -44.444444444444
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F16").Value = 16
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 42
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("H16").Value = -61.904761904761
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I16").Value = 40
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 70
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("K16").Value = -42.857142857142
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = -27.272727272727
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M16").Value = -39.393939393939
$ws.Range("M16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N16").Value = -84.615384615384
$ws.Range("N16").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 15
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 9
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = 66.666666666666
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 43
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 44
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("H17").Value = -2.272727272727
$ws.Range("H17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I17").Value = 88
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 86
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("K17").Value = 2.325581395348
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = -6.382978723404
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M17").Value = 175
$ws.Range("M17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N17").Value = 60
$ws.Range("N17").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 18 (Burglary)
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 4
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = -75
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G18").Value = 26
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("H18").Value = -38.461538461538
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I18").Value = 28
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("J18").Value = 51
$ws.Range("J18").NumberFormat = "#,##0"
$ws.Range("K18").Value = -45.098039215686
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = -22.222222222222
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M18").Value = -57.575757575757
$ws.Range("M18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N18").Value = -92.893401015228
$ws.Range("N18").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 10
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 27
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("E19").Value = -62.962962962963
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F19").Value = 45
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 92
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("H19").Value = -51.086956521739
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I19").Value = 102
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 170
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("K19").Value = -40
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = -10.526315789473
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M19").Value = 12.087912087912
$ws.Range("M19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N19").Value = -40
$ws.Range("N19").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 5
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 5
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("F20").Value = 18
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("G20").Value = 17
$ws.Range("G20").NumberFormat = "#,##0"
$ws.Range("H20").Value = 5.882352941176
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I20").Value = 28
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("J20").Value = 39
$ws.Range("J20").NumberFormat = "#,##0"
$ws.Range("K20").Value = -28.205128205128
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = -31.707317073170
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M20").Value = 27.272727272727
$ws.Range("M20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N20").Value = -91.715976331360
$ws.Range("N20").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 21 (TOTAL, bold style 17/18)
$ws.Range("C21").Value = 36
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 54
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("E21").Value = -33.333333333333
$ws.Range("E21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("F21").Value = 140
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 223
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("H21").Value = -37.219730941704
$ws.Range("H21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("I21").Value = 290
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 423
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("K21").Value = -31.442080378250
$ws.Range("K21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("L21").Value = -15.942028985507
$ws.Range("L21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("M21").Value = 2.836879432624
$ws.Range("M21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("N21").Value = -76.326530612244
$ws.Range("N21").NumberFormat = "#,##0.00;""-""#,##0.00"

# Row 22 (Transit): C becomes a number, D/E become the "no data" text
# markers (the reverse of row 15's change).
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 2
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("H22").Value = -33.333333333333
$ws.Range("H22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I22").Value = 7
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("K22").Value = -22.222222222222
$ws.Range("K22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L22").Value = 133.333333333333
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M22").Value = 16.666666666666
$ws.Range("M22").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 59
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 58
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("E24").Value = 1.724137931034
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F24").Value = 204
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 261
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("H24").Value = -21.839080459770
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I24").Value = 383
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 496
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("K24").Value = -22.782258064516
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = -15.265486725663
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M24").Value = 86.829268292682
$ws.Range("M24").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 50
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 57
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("E25").Value = -12.280701754386
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F25").Value = 166
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 233
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("H25").Value = -28.755364806867
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I25").Value = 304
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 409
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("K25").Value = -25.672371638141
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = -8.708708708708
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 15
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 28
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -46.428571428571
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 58
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 94
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = -38.297872340425
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 120
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 173
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = -30.635838150289
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L26").Value = -12.408759124087
$ws.Range("L26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M26").Value = 55.844155844155
$ws.Range("M26").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 27 (UCR Rape*): C/D/E switch to the "no data" text markers.
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("D14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("G27").Value = 2
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").Value = 50
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L27").Value = -16.666666666666
$ws.Range("L27").NumberFormat = "#,##0.0;""-""#,##0.0"

# Row 28 (Other Sex Crimes): C becomes a number.
$ws.Range("C28").Value = 4
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 300
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G28").Value = 13
$ws.Range("G28").NumberFormat = "#,##0"
$ws.Range("H28").Value = -15.384615384615
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I28").Value = 18
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 20
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("K28").Value = -10
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = -18.181818181818
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
